$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new worksheet row above row 37 (shifts existing rows 37-41 down to 38-42),
# then grow the table to cover the newly-inserted row.
$ws.Rows.Item(37).Insert()
$lo.Resize($ws.Range("A1:D42"))

# Fill in the new hike entry: Stevens Creek Trail, 1.2 mi, 500 ft, easy.
$ws.Cells.Item(37, 1).Value = "Stevens Creek Trail"
$ws.Cells.Item(37, 2).Value = 1.2
$ws.Cells.Item(37, 3).Value = 500
$ws.Cells.Item(37, 4).Value = "easy"

# Match the post-edit selection left behind by Excel (one row past the new last row).
$ws.Range("D43").Select() | Out-Null
